# Auto-generated edit script: updates numeric cells per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
  # Row 129
  $ws.Cells.Item(129, 8).Value = 43210692
  $ws.Cells.Item(129, 9).Value = 111111590
  $ws.Cells.Item(129, 10).Value = 2470152.5
  $ws.Cells.Item(129, 11).Value = 333334770
  $ws.Cells.Item(129, 12).Value = 7410457.5
  $ws.Cells.Item(129, 13).Value = -333329770
  $ws.Cells.Item(129, 14).Value = -7420457.5
  # Row 132
  $ws.Cells.Item(132, 8).Value = 4789.5
  $ws.Cells.Item(132, 9).Value = 6213.5713
  $ws.Cells.Item(132, 11).Value = 18640.7139
  $ws.Cells.Item(132, 13).Value = -16110.7139
  # Row 138
  $ws.Cells.Item(138, 8).Value = 3968.2524
  $ws.Cells.Item(138, 9).Value = 1136.0476
  $ws.Cells.Item(138, 10).Value = 4730.769
  $ws.Cells.Item(138, 11).Value = 3408.142800000001
  $ws.Cells.Item(138, 12).Value = 14192.307
  $ws.Cells.Item(138, 13).Value = 1731.857199999999
  $ws.Cells.Item(138, 14).Value = -24472.307

$ws = $wb.Worksheets.Item("ARM")
  # Row 2
  $ws.Cells.Item(2, 8).Value = 48055.316
  $ws.Cells.Item(2, 9).Value = 60630.59
  $ws.Cells.Item(2, 11).Value = 60630.59
  $ws.Cells.Item(2, 13).Value = -60517.59
  # Row 32
  $ws.Cells.Item(32, 8).Value = 15618.408
  $ws.Cells.Item(32, 9).Value = 12805.688
  $ws.Cells.Item(32, 11).Value = 12805.688
  $ws.Cells.Item(32, 13).Value = -12518.688
  # Row 45
  $ws.Cells.Item(45, 8).Value = 6519.6523
  $ws.Cells.Item(45, 9).Value = 8574.143
  $ws.Cells.Item(45, 11).Value = 8574.143
  $ws.Cells.Item(45, 13).Value = -8197.143
  # Row 116
  $ws.Cells.Item(116, 8).Value = 48055.316
  $ws.Cells.Item(116, 9).Value = 60630.59
  $ws.Cells.Item(116, 11).Value = 60630.59
  $ws.Cells.Item(116, 13).Value = -58336.59

$ws = $wb.Worksheets.Item("BSM")
  # Row 3
  $ws.Cells.Item(3, 8).Value = 48055.316
  $ws.Cells.Item(3, 9).Value = 60630.59
  $ws.Cells.Item(3, 11).Value = 60630.59
  $ws.Cells.Item(3, 13).Value = -60516.59
  # Row 58
  $ws.Cells.Item(58, 8).Value = 25650
  $ws.Cells.Item(58, 10).Value = 25650
  $ws.Cells.Item(58, 12).Value = 25650
  $ws.Cells.Item(58, 14).Value = -26238
  # Row 86
  $ws.Cells.Item(86, 8).Value = 2556
  $ws.Cells.Item(86, 9).Value = 2248.1177
  $ws.Cells.Item(86, 10).Value = 3602.8
  $ws.Cells.Item(86, 11).Value = 2248.1177
  $ws.Cells.Item(86, 12).Value = 3602.8
  $ws.Cells.Item(86, 13).Value = -1125.1177
  $ws.Cells.Item(86, 14).Value = -5848.8
  # Row 89
  $ws.Cells.Item(89, 8).Value = 2556
  $ws.Cells.Item(89, 9).Value = 2248.1177
  $ws.Cells.Item(89, 10).Value = 3602.8
  $ws.Cells.Item(89, 11).Value = 11240.5885
  $ws.Cells.Item(89, 12).Value = 18014
  $ws.Cells.Item(89, 13).Value = -5624.588499999998
  $ws.Cells.Item(89, 14).Value = -29246
  # Row 94
  $ws.Cells.Item(94, 8).Value = 1662.0769
  $ws.Cells.Item(94, 9).Value = 1383.9166
  $ws.Cells.Item(94, 10).Value = 5000
  $ws.Cells.Item(94, 11).Value = 1383.9166
  $ws.Cells.Item(94, 12).Value = 5000
  $ws.Cells.Item(94, 13).Value = -932.9166
  $ws.Cells.Item(94, 14).Value = -5902
  # Row 99
  $ws.Cells.Item(99, 8).Value = 40001908
  $ws.Cells.Item(99, 9).Value = 47620750
  $ws.Cells.Item(99, 11).Value = 47620750
  $ws.Cells.Item(99, 13).Value = -47619252
  # Row 105
  $ws.Cells.Item(105, 8).Value = 3154.2104
  $ws.Cells.Item(105, 9).Value = 3271.5386
  $ws.Cells.Item(105, 11).Value = 3271.5386
  $ws.Cells.Item(105, 13).Value = -1524.5386
  # Row 134
  $ws.Cells.Item(134, 8).Value = 2287.739
  $ws.Cells.Item(134, 9).Value = 2537.6365
  $ws.Cells.Item(134, 10).Value = 2058.6667
  $ws.Cells.Item(134, 11).Value = 7612.9095
  $ws.Cells.Item(134, 12).Value = 6176.000100000001
  $ws.Cells.Item(134, 13).Value = -5077.9095
  $ws.Cells.Item(134, 14).Value = -11246.0001

$ws = $wb.Worksheets.Item("CRP")
  # Row 31
  $ws.Cells.Item(31, 8).Value = 2674.3076
  $ws.Cells.Item(31, 9).Value = 2443.5
  $ws.Cells.Item(31, 11).Value = 2443.5
  $ws.Cells.Item(31, 13).Value = -2148.5
  # Row 34
  $ws.Cells.Item(34, 8).Value = 2674.3076
  $ws.Cells.Item(34, 9).Value = 2443.5
  $ws.Cells.Item(34, 11).Value = 2443.5
  $ws.Cells.Item(34, 13).Value = -2241.5
  # Row 134
  $ws.Cells.Item(134, 8).Value = 4887.4287
  $ws.Cells.Item(134, 9).Value = 5178
  $ws.Cells.Item(134, 11).Value = 15534
  $ws.Cells.Item(134, 13).Value = -12999

$ws = $wb.Worksheets.Item("CUL")
  # Row 122
  $ws.Cells.Item(122, 8).Value = 360.6316
  $ws.Cells.Item(122, 9).Value = 360.6316
  $ws.Cells.Item(122, 10).Value = 0
  $ws.Cells.Item(122, 11).Value = 3245.6844
  $ws.Cells.Item(122, 12).Value = 0
  $ws.Cells.Item(122, 13).Value = -795.6844000000001
  $ws.Cells.Item(122, 14).ClearContents()
  # Row 136
  $ws.Cells.Item(136, 8).Value = 3807.077
  $ws.Cells.Item(136, 9).Value = 2922.7144
  $ws.Cells.Item(136, 10).Value = 4838.8335
  $ws.Cells.Item(136, 11).Value = 8768.143199999999
  $ws.Cells.Item(136, 12).Value = 14516.5005
  $ws.Cells.Item(136, 13).Value = -3668.143199999999
  $ws.Cells.Item(136, 14).Value = -24716.5005

$ws = $wb.Worksheets.Item("GSM")
  # Row 132
  $ws.Cells.Item(132, 8).Value = 4028.3333
  $ws.Cells.Item(132, 9).Value = 3750
  $ws.Cells.Item(132, 10).Value = 4084
  $ws.Cells.Item(132, 11).Value = 11250
  $ws.Cells.Item(132, 12).Value = 12252
  $ws.Cells.Item(132, 13).Value = -8720
  $ws.Cells.Item(132, 14).Value = -17312

$ws = $wb.Worksheets.Item("LTW")
  # Row 61
  $ws.Cells.Item(61, 8).Value = 0
  $ws.Cells.Item(61, 9).Value = 0
  $ws.Cells.Item(61, 10).Value = 0
  $ws.Cells.Item(61, 11).Value = 0
  $ws.Cells.Item(61, 12).Value = 0
  $ws.Cells.Item(61, 13).ClearContents()
  $ws.Cells.Item(61, 14).ClearContents()
  # Row 68
  $ws.Cells.Item(68, 8).Value = 2283.3333
  $ws.Cells.Item(68, 9).Value = 2225
  $ws.Cells.Item(68, 10).Value = 2400
  $ws.Cells.Item(68, 11).Value = 2225
  $ws.Cells.Item(68, 12).Value = 2400
  $ws.Cells.Item(68, 13).Value = -1476
  $ws.Cells.Item(68, 14).Value = -3898
  # Row 71
  $ws.Cells.Item(71, 8).Value = 2283.3333
  $ws.Cells.Item(71, 9).Value = 2225
  $ws.Cells.Item(71, 10).Value = 2400
  $ws.Cells.Item(71, 11).Value = 11125
  $ws.Cells.Item(71, 12).Value = 12000
  $ws.Cells.Item(71, 13).Value = -7381
  $ws.Cells.Item(71, 14).Value = -19488
  # Row 93
  $ws.Cells.Item(93, 8).Value = 1107.0714
  $ws.Cells.Item(93, 9).Value = 982.9167
  $ws.Cells.Item(93, 11).Value = 982.9167
  $ws.Cells.Item(93, 13).Value = 265.0833
  # Row 105
  $ws.Cells.Item(105, 8).Value = 0
  $ws.Cells.Item(105, 10).Value = 0
  $ws.Cells.Item(105, 12).Value = 0
  $ws.Cells.Item(105, 14).ClearContents()
  # Row 113
  $ws.Cells.Item(113, 8).Value = 0
  $ws.Cells.Item(113, 9).Value = 0
  $ws.Cells.Item(113, 10).Value = 0
  $ws.Cells.Item(113, 11).Value = 0
  $ws.Cells.Item(113, 12).Value = 0
  $ws.Cells.Item(113, 13).ClearContents()
  $ws.Cells.Item(113, 14).ClearContents()
  # Row 125
  $ws.Cells.Item(125, 8).Value = 37326.668
  $ws.Cells.Item(125, 10).Value = 37326.668
  $ws.Cells.Item(125, 12).Value = 37326.668
  $ws.Cells.Item(125, 14).Value = -47166.668
  # Row 127
  $ws.Cells.Item(127, 8).Value = 0
  $ws.Cells.Item(127, 10).Value = 0
  $ws.Cells.Item(127, 12).Value = 0
  $ws.Cells.Item(127, 14).Value = 0
  # Row 136
  $ws.Cells.Item(136, 8).Value = 4129.625
  $ws.Cells.Item(136, 9).Value = 1812.3334
  $ws.Cells.Item(136, 10).Value = 5520
  $ws.Cells.Item(136, 11).Value = 5437.0002
  $ws.Cells.Item(136, 12).Value = 16560
  $ws.Cells.Item(136, 13).Value = -2887.0002
  $ws.Cells.Item(136, 14).Value = -21660
  # Row 139
  $ws.Cells.Item(139, 8).Value = 43287.855
  $ws.Cells.Item(139, 10).Value = 43287.855
  $ws.Cells.Item(139, 12).Value = 43287.855
  $ws.Cells.Item(139, 14).Value = -53567.855
